# Updated cryptos list — refresh Price (D) and Volume(1h) (E) columns with
# latest scraped figures. A handful of Price cells are plain numeric-looking
# strings (e.g. "217.50"); Excel's Range.Value setter auto-coerces those to
# numbers, which would lose the trailing zero / change the stored type. For
# those cells we briefly force a Text number format, assign the value, then
# restore the cell to the "Normal" style so no stray formatting is left
# behind.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($addr, $value) {
    $cell = $ws.Range($addr)
    $cell.NumberFormat = "@"
    $cell.Value = $value
    $cell.Style = "Normal"
}

$ws.Range("D2").Value = "26.713.50"
$ws.Range("D3").Value = "1.637.16"
$ws.Range("E3").Value = "  -0.70%  "
Set-TextValue "D5" "217.50"
$ws.Range("E5").Value = "  +0.54%  "
Set-TextValue "D6" "0.502"
$ws.Range("E6").Value = "  -0.87%  "
$ws.Range("E7").Value = "  -0.02%  "
$ws.Range("E8").Value = "  -0.69%  "
$ws.Range("E9").Value = "  -0.80%  "
Set-TextValue "D10" "19.08"
$ws.Range("E10").Value = "  -0.56%  "
Set-TextValue "D11" "0.0845"
$ws.Range("E11").Value = "  +0.32%  "
$ws.Range("D12").Value = "1.863.14"
$ws.Range("D13").Value = "1.635.30"
$ws.Range("E13").Value = "  -1.12%  "
$ws.Range("E14").Value = "  -1.21%  "
$ws.Range("E15").Value = "  -1.41%  "
Set-TextValue "D16" "64.44"
$ws.Range("E16").Value = "  -1.58%  "
$ws.Range("D17").Value = "26.695.66"
$ws.Range("E17").Value = "  -0.31%  "
$ws.Range("D18").Value = "0.0₃0726"
$ws.Range("E18").Value = "  -2.45%  "
Set-TextValue "D19" "211.30"
$ws.Range("E19").Value = "  -3.37%  "
$ws.Range("E20").Value = "  +0.01%  "
Set-TextValue "D21" "4.33"
$ws.Range("E21").Value = "  -0.89%  "
$ws.Range("E22").Value = "  -1.39%  "
Set-TextValue "D23" "2.30"
$ws.Range("E23").Value = "  -2.51%  "
Set-TextValue "D24" "9.26"
$ws.Range("E24").Value = "  -2.67%  "
Set-TextValue "D25" "146.42"
$ws.Range("E25").Value = "  -0.12%  "
$ws.Range("E26").Value = "  +0.00%  "
Set-TextValue "D27" "0.117"
$ws.Range("E27").Value = "  -2.21%  "
Set-TextValue "D28" "7.08"
$ws.Range("E28").Value = "  -0.56%  "
Set-TextValue "D29" "15.55"
$ws.Range("E29").Value = "  -1.35%  "
$ws.Range("E30").Value = "  -2.46%  "
$ws.Range("E31").Value = "  +0.59%  "
$ws.Range("E32").Value = "  -0.24%  "
$ws.Range("E33").Value = "  -1.39%  "
$ws.Range("D34").Value = "1.272.17"
$ws.Range("E34").Value = "  -0.65%  "
$ws.Range("E35").Value = "  -1.31%  "
Set-TextValue "D36" "2.44"
$ws.Range("E36").Value = "  +0.08%  "
$ws.Range("E37").Value = "  -1.91%  "
$ws.Range("E38").Value = "  -1.87%  "
$ws.Range("E39").Value = "  -2.72%  "
$ws.Range("E40").Value = "  -0.03%  "
$ws.Range("E41").Value = "  -1.55%  "
$ws.Range("E42").Value = "  -2.53%  "
$ws.Range("E43").Value = "  -3.64%  "
$ws.Range("D44").Value = "1.774.29"
$ws.Range("E44").Value = "  -0.86%  "
Set-TextValue "D45" "91.41"
$ws.Range("E45").Value = "  -0.67%  "
Set-TextValue "D46" "60.31"
$ws.Range("E46").Value = "  +0.93%  "
$ws.Range("E47").Value = "  -1.55%  "
$ws.Range("E48").Value = "  +0.42%  "
Set-TextValue "D49" "7.55"
$ws.Range("E49").Value = "  -2.63%  "
Set-TextValue "D50" "0.0961"
$ws.Range("E50").Value = "  -0.90%  "
$ws.Range("B51").Value = "Mantle"
$ws.Range("C51").Value = "https://coinranking.com/coin/BoI4ux0nd+mantle-mnt"
Set-TextValue "D51" "0.406"
$ws.Range("E51").Value = "  -0.35%  "
